$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 824.7037
$ws.Range("J129").Value = 898.63635
$ws.Range("L129").Value = 2695.90905
$ws.Range("N129").Value = -12695.90905
$ws.Range("H135").Value = 31260224
$ws.Range("I135").Value = 1325.8182
$ws.Range("K135").Value = 11932.3638
$ws.Range("M135").Value = -9397.363799999999
$ws.Range("H137").Value = 40019.406
$ws.Range("I137").Value = 3869.2307
$ws.Range("J137").Value = 73587.42999999999
$ws.Range("K137").Value = 11607.6921
$ws.Range("L137").Value = 220762.29
$ws.Range("M137").Value = -9057.6921
$ws.Range("N137").Value = -225862.29
$ws.Range("H138").Value = 2468.4707
$ws.Range("I138").Value = 920.25
$ws.Range("J138").Value = 3312.9546
$ws.Range("K138").Value = 2760.75
$ws.Range("L138").Value = 9938.863799999999
$ws.Range("M138").Value = 2379.25
$ws.Range("N138").Value = -20218.8638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24953.195
$ws.Range("I32").Value = 31049.371
$ws.Range("K32").Value = 31049.371
$ws.Range("M32").Value = -30762.371
$ws.Range("H61").Value = 465021.94
$ws.Range("I61").Value = 752847.5600000001
$ws.Range("K61").Value = 752847.5600000001
$ws.Range("M61").Value = -752635.5600000001
$ws.Range("H76").Value = 26000
$ws.Range("J76").Value = 26000
$ws.Range("L76").Value = 26000
$ws.Range("N76").Value = -26676
$ws.Range("H79").Value = 26000
$ws.Range("J79").Value = 26000
$ws.Range("L79").Value = 26000
$ws.Range("N79").Value = -28340
$ws.Range("H110").Value = 3305.3635
$ws.Range("I110").Value = 2297.375
$ws.Range("J110").Value = 5993.3335
$ws.Range("K110").Value = 2297.375
$ws.Range("L110").Value = 5993.3335
$ws.Range("M110").Value = -252.375
$ws.Range("N110").Value = -10083.3335
$ws.Range("H132").Value = 31457.295
$ws.Range("I132").Value = 2162.4285
$ws.Range("J132").Value = 168166.67
$ws.Range("K132").Value = 6487.2855
$ws.Range("L132").Value = 504500.01
$ws.Range("M132").Value = -3957.2855
$ws.Range("N132").Value = -509560.01
$ws.Range("H135").Value = 22101.9
$ws.Range("J135").Value = 22101.9
$ws.Range("L135").Value = 22101.9
$ws.Range("N135").Value = -32241.9
$ws.Range("H136").Value = 465021.94
$ws.Range("I136").Value = 752847.5600000001
$ws.Range("K136").Value = 2258542.68
$ws.Range("M136").Value = -2255992.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 53916.25
$ws.Range("I134").Value = 59779.168
$ws.Range("K134").Value = 179337.504
$ws.Range("M134").Value = -176802.504

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11618.137
$ws.Range("I31").Value = 16055.357
$ws.Range("J31").Value = 3853
$ws.Range("K31").Value = 16055.357
$ws.Range("L31").Value = 3853
$ws.Range("M31").Value = -15760.357
$ws.Range("N31").Value = -4443
$ws.Range("H34").Value = 11618.137
$ws.Range("I34").Value = 16055.357
$ws.Range("J34").Value = 3853
$ws.Range("K34").Value = 16055.357
$ws.Range("L34").Value = 3853
$ws.Range("M34").Value = -15853.357
$ws.Range("N34").Value = -4257
$ws.Range("H58").Value = 30744.705
$ws.Range("I58").Value = 1310.6666
$ws.Range("J58").Value = 251500
$ws.Range("K58").Value = 1310.6666
$ws.Range("L58").Value = 251500
$ws.Range("M58").Value = -1107.6666
$ws.Range("N58").Value = -251906
$ws.Range("H122").Value = 2215.75
$ws.Range("I122").Value = 2802.4
$ws.Range("J122").Value = 1238
$ws.Range("K122").Value = 8407.200000000001
$ws.Range("L122").Value = 3714
$ws.Range("M122").Value = -5957.200000000001
$ws.Range("N122").Value = -8614
$ws.Range("H136").Value = 30744.705
$ws.Range("I136").Value = 1310.6666
$ws.Range("J136").Value = 251500
$ws.Range("K136").Value = 3931.9998
$ws.Range("L136").Value = 754500
$ws.Range("M136").Value = -1381.9998
$ws.Range("N136").Value = -759600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2412.2222
$ws.Range("I51").Value = 1500
$ws.Range("J51").Value = 2672.8572
$ws.Range("K51").Value = 4500
$ws.Range("L51").Value = 8018.571599999999
$ws.Range("M51").Value = -4040
$ws.Range("N51").Value = -8938.571599999999
$ws.Range("H56").Value = 6107.0713
$ws.Range("I56").Value = 6107.0713
$ws.Range("K56").Value = 6107.0713
$ws.Range("M56").Value = -5577.0713
$ws.Range("H68").Value = 3890.4722
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 3958.7715
$ws.Range("K68").Value = 4500
$ws.Range("L68").Value = 11876.3145
$ws.Range("M68").Value = -3689
$ws.Range("N68").Value = -13498.3145
$ws.Range("H71").Value = 3890.4722
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 3958.7715
$ws.Range("K71").Value = 13500
$ws.Range("L71").Value = 35628.9435
$ws.Range("M71").Value = -9444
$ws.Range("N71").Value = -43740.9435
$ws.Range("H107").Value = 4069.0938
$ws.Range("I107").Value = 25449
$ws.Range("J107").Value = 1014.8214
$ws.Range("K107").Value = 76347
$ws.Range("L107").Value = 3044.4642
$ws.Range("M107").Value = -74427
$ws.Range("N107").Value = -6884.4642
$ws.Range("H109").Value = 4464.4287
$ws.Range("I109").Value = 1102
$ws.Range("J109").Value = 6332.4443
$ws.Range("K109").Value = 3306
$ws.Range("L109").Value = 18997.3329
$ws.Range("M109").Value = -2266
$ws.Range("N109").Value = -21077.3329
$ws.Range("H113").Value = 20618.2
$ws.Range("I113").Value = 50350.5
$ws.Range("J113").Value = 796.6667
$ws.Range("K113").Value = 151051.5
$ws.Range("L113").Value = 2390.0001
$ws.Range("M113").Value = -148881.5
$ws.Range("N113").Value = -6730.0001
$ws.Range("H131").Value = 135979.78
$ws.Range("J131").Value = 150103.64
$ws.Range("L131").Value = 450310.92
$ws.Range("N131").Value = -460390.92
$ws.Range("H132").Value = 740.875
$ws.Range("I132").Value = 704.9091
$ws.Range("J132").Value = 820
$ws.Range("K132").Value = 6344.1819
$ws.Range("L132").Value = 7380
$ws.Range("M132").Value = -3814.1819
$ws.Range("N132").Value = -12440
$ws.Range("H140").Value = 1856.6471
$ws.Range("I140").Value = 1473.9286
$ws.Range("K140").Value = 4421.7858
$ws.Range("M140").Value = 758.2142000000003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 8628
$ws.Range("I102").Value = 10670.667
$ws.Range("K102").Value = 10670.667
$ws.Range("M102").Value = -9048.666999999999
$ws.Range("H112").Value = 20000
$ws.Range("J112").Value = 20000
$ws.Range("L112").Value = 20000
$ws.Range("N112").Value = -22216
$ws.Range("H122").Value = 3869.4167
$ws.Range("I122").Value = 2938.111
$ws.Range("J122").Value = 6663.3335
$ws.Range("K122").Value = 8814.332999999999
$ws.Range("L122").Value = 19990.0005
$ws.Range("M122").Value = -6364.332999999999
$ws.Range("N122").Value = -24890.0005
$ws.Range("H126").Value = 5484.8623
$ws.Range("I126").Value = 4919
$ws.Range("K126").Value = 14757
$ws.Range("M126").Value = -12287
$ws.Range("H132").Value = 87410.61
$ws.Range("I132").Value = 95581.91
$ws.Range("J132").Value = 74570
$ws.Range("K132").Value = 286745.73
$ws.Range("L132").Value = 223710
$ws.Range("M132").Value = -284215.73
$ws.Range("N132").Value = -228770

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2080.3333
$ws.Range("I22").Value = 2320
$ws.Range("J22").Value = 882
$ws.Range("K22").Value = 2320
$ws.Range("L22").Value = 882
$ws.Range("M22").Value = -2025
$ws.Range("N22").Value = -1472
$ws.Range("H27").Value = 2080.3333
$ws.Range("I27").Value = 2320
$ws.Range("J27").Value = 882
$ws.Range("K27").Value = 2320
$ws.Range("L27").Value = 882
$ws.Range("M27").Value = -2213
$ws.Range("N27").Value = -1096
$ws.Range("H40").Value = 89754.69500000001
$ws.Range("I40").Value = 104619.18
$ws.Range("K40").Value = 104619.18
$ws.Range("M40").Value = -104483.18
$ws.Range("H100").Value = 4314.2856
$ws.Range("I100").Value = 1850
$ws.Range("K100").Value = 1850
$ws.Range("M100").Value = -1309

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H122").Value = 1562.28
$ws.Range("I122").Value = 1410.1428
$ws.Range("J122").Value = 2361
$ws.Range("K122").Value = 4230.428400000001
$ws.Range("L122").Value = 7083
$ws.Range("M122").Value = -1780.428400000001
$ws.Range("N122").Value = -11983
$ws.Range("H132").Value = 3588.6667
$ws.Range("J132").Value = 4499.3335
$ws.Range("L132").Value = 13498.0005
$ws.Range("N132").Value = -18558.0005
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
